$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "41-14="
$t.Cell(1,2).Range.Text = "81-76="
$t.Cell(1,3).Range.Text = "51-12="
$t.Cell(1,4).Range.Text = "51-16="
$t.Cell(1,5).Range.Text = "53+38="
$t.Cell(2,1).Range.Text = "44-17="
$t.Cell(2,2).Range.Text = "54+7="
$t.Cell(2,3).Range.Text = "12+59="
$t.Cell(2,4).Range.Text = "28+46="
$t.Cell(2,5).Range.Text = "75+6="
$t.Cell(3,1).Range.Text = "61-17="
$t.Cell(3,2).Range.Text = "70-31="
$t.Cell(3,3).Range.Text = "26+69="
$t.Cell(3,4).Range.Text = "51-24="
$t.Cell(3,5).Range.Text = "49+19="
$t.Cell(4,1).Range.Text = "58+39="
$t.Cell(4,2).Range.Text = "65-59="
$t.Cell(4,3).Range.Text = "39+26="
$t.Cell(4,4).Range.Text = "72-56="
$t.Cell(4,5).Range.Text = "16+66="
$t.Cell(5,1).Range.Text = "53-7="
$t.Cell(5,2).Range.Text = "90-25="
$t.Cell(5,3).Range.Text = "83-27="
$t.Cell(5,4).Range.Text = "56+25="
$t.Cell(5,5).Range.Text = "57+5="
$t.Cell(6,1).Range.Text = "17+7="
$t.Cell(6,2).Range.Text = "94-68="
$t.Cell(6,3).Range.Text = "31-4="
$t.Cell(6,4).Range.Text = "52-15="
$t.Cell(6,5).Range.Text = "55+28="
$t.Cell(7,1).Range.Text = "62-45="
$t.Cell(7,2).Range.Text = "5+28="
$t.Cell(7,3).Range.Text = "62-36="
$t.Cell(7,4).Range.Text = "82-29="
$t.Cell(7,5).Range.Text = "17+44="
$t.Cell(8,1).Range.Text = "9+7="
$t.Cell(8,2).Range.Text = "80-33="
$t.Cell(8,3).Range.Text = "85-26="
$t.Cell(8,4).Range.Text = "5+36="
$t.Cell(8,5).Range.Text = "81-15="
$t.Cell(9,1).Range.Text = "39+44="
$t.Cell(9,2).Range.Text = "71-47="
$t.Cell(9,3).Range.Text = "14-5="
$t.Cell(9,4).Range.Text = "8+4="
$t.Cell(9,5).Range.Text = "70-55="
$t.Cell(10,1).Range.Text = "71-65="
$t.Cell(10,2).Range.Text = "8+63="
$t.Cell(10,3).Range.Text = "85-18="
$t.Cell(10,4).Range.Text = "93-67="
$t.Cell(10,5).Range.Text = "66+26="
$t.Cell(11,1).Range.Text = "44-9="
$t.Cell(11,2).Range.Text = "91-62="
$t.Cell(11,3).Range.Text = "73-64="
$t.Cell(11,4).Range.Text = "26+5="
$t.Cell(11,5).Range.Text = "27+39="
$t.Cell(12,1).Range.Text = "49+49="
$t.Cell(12,2).Range.Text = "81-49="
$t.Cell(12,3).Range.Text = "97-48="
$t.Cell(12,4).Range.Text = "30-1="
$t.Cell(12,5).Range.Text = "58+24="
$t.Cell(13,1).Range.Text = "70-47="
$t.Cell(13,2).Range.Text = "48-19="
$t.Cell(13,3).Range.Text = "77-68="
$t.Cell(13,4).Range.Text = "52+19="
$t.Cell(13,5).Range.Text = "81-45="
$t.Cell(14,1).Range.Text = "62-13="
$t.Cell(14,2).Range.Text = "14+27="
$t.Cell(14,3).Range.Text = "91-3="
$t.Cell(14,4).Range.Text = "36+46="
$t.Cell(14,5).Range.Text = "73-29="
$t.Cell(15,1).Range.Text = "27+17="
$t.Cell(15,2).Range.Text = "7+69="
$t.Cell(15,3).Range.Text = "33-9="
$t.Cell(15,4).Range.Text = "84-39="
$t.Cell(15,5).Range.Text = "54+37="
$t.Cell(16,1).Range.Text = "14+77="
$t.Cell(16,2).Range.Text = "50-31="
$t.Cell(16,3).Range.Text = "23+59="
$t.Cell(16,4).Range.Text = "58+34="
$t.Cell(16,5).Range.Text = "8+35="
$t.Cell(17,1).Range.Text = "59+23="
$t.Cell(17,2).Range.Text = "7+36="
$t.Cell(17,3).Range.Text = "75+18="
$t.Cell(17,4).Range.Text = "18+75="
$t.Cell(17,5).Range.Text = "36+58="
$t.Cell(18,1).Range.Text = "50-44="
$t.Cell(18,2).Range.Text = "28+26="
$t.Cell(18,3).Range.Text = "18+14="
$t.Cell(18,4).Range.Text = "55-16="
$t.Cell(18,5).Range.Text = "40-16="
$t.Cell(19,1).Range.Text = "29+9="
$t.Cell(19,2).Range.Text = "81-23="
$t.Cell(19,3).Range.Text = "82-6="
$t.Cell(19,4).Range.Text = "59+33="
$t.Cell(19,5).Range.Text = "95-36="
$t.Cell(20,1).Range.Text = "98-59="
$t.Cell(20,2).Range.Text = "78+19="
$t.Cell(20,3).Range.Text = "5+46="
$t.Cell(20,4).Range.Text = "81-67="
$t.Cell(20,5).Range.Text = "29+66="
